$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Apply formatting to the new rows (639:654) ---
# Copy the full row format from row 638 (an existing populated data row) as the
# baseline for the newly appended rows 639:654.
$ws.Range("A638:I638").Copy()
$ws.Range("A639:I654").PasteSpecial(-4122)

# Fix up the "G" cells that should stay empty (no injury location) to use the
# centered blank style, matching row 637's empty G cell format.
$ws.Range("G637").Copy()
$ws.Range("G639").PasteSpecial(-4122)
$ws.Range("G645").PasteSpecial(-4122)
$ws.Range("G646").PasteSpecial(-4122)
$ws.Range("G648").PasteSpecial(-4122)
$ws.Range("G651").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 2: Enter the new row values (training log entries for 2025-12-03) ---
$ws.Range("A639").Value = 45994
$ws.Range("B639").Value = "Karim Belmahi"
$ws.Range("C639").Value = 70
$ws.Range("D639").Value = 8
$ws.Range("E639").Value = 7
$ws.Range("F639").Value = 0
$ws.Range("H639").Value = 10
$ws.Range("A640").Value = 45994
$ws.Range("B640").Value = "Yoann Martelat"
$ws.Range("C640").Value = 70
$ws.Range("D640").Value = 7
$ws.Range("E640").Value = 9
$ws.Range("F640").Value = 3
$ws.Range("G640").Value = "Genou"
$ws.Range("H640").Value = 3
$ws.Range("A641").Value = 45994
$ws.Range("B641").Value = "Amir Etien"
$ws.Range("C641").Value = 70
$ws.Range("D641").Value = 6
$ws.Range("E641").Value = 7
$ws.Range("F641").Value = 4
$ws.Range("G641").Value = "Ischio droit"
$ws.Range("H641").Value = 2
$ws.Range("A642").Value = 45994
$ws.Range("B642").Value = "Jeremie Laurent"
$ws.Range("C642").Value = 70
$ws.Range("D642").Value = 5
$ws.Range("E642").Value = 5
$ws.Range("F642").Value = 2
$ws.Range("G642").Value = "Ischio"
$ws.Range("H642").Value = 6
$ws.Range("A643").Value = 45994
$ws.Range("B643").Value = "Yoan Zouma"
$ws.Range("C643").Value = 70
$ws.Range("D643").Value = 6
$ws.Range("E643").Value = 8
$ws.Range("F643").Value = 5
$ws.Range("G643").Value = "Ischio"
$ws.Range("H643").Value = 3
$ws.Range("A644").Value = 45994
$ws.Range("B644").Value = "Levy Ndoutoume"
$ws.Range("C644").Value = 70
$ws.Range("D644").Value = 6
$ws.Range("E644").Value = 7
$ws.Range("F644").Value = 1
$ws.Range("G644").Value = "Ischio"
$ws.Range("H644").Value = 3
$ws.Range("A645").Value = 45994
$ws.Range("B645").Value = "Emmanuel Valey"
$ws.Range("C645").Value = 70
$ws.Range("D645").Value = 8
$ws.Range("E645").Value = 6
$ws.Range("F645").Value = 0
$ws.Range("H645").Value = 6
$ws.Range("A646").Value = 45994
$ws.Range("B646").Value = "Ilan Ihaddadene"
$ws.Range("C646").Value = 70
$ws.Range("D646").Value = 6
$ws.Range("E646").Value = 7
$ws.Range("F646").Value = 0
$ws.Range("H646").Value = 5
$ws.Range("A647").Value = 45994
$ws.Range("B647").Value = "Karahali Souaré"
$ws.Range("C647").Value = 70
$ws.Range("D647").Value = 3
$ws.Range("E647").Value = 3
$ws.Range("F647").Value = 5
$ws.Range("G647").Value = "Cheville"
$ws.Range("H647").Value = 3
$ws.Range("A648").Value = 45994
$ws.Range("B648").Value = "Mattheo Haon"
$ws.Range("C648").Value = 70
$ws.Range("D648").Value = 6
$ws.Range("E648").Value = 5
$ws.Range("F648").Value = 0
$ws.Range("H648").Value = 10
$ws.Range("A649").Value = 45994
$ws.Range("B649").Value = "Romain Thunet"
$ws.Range("C649").Value = 70
$ws.Range("D649").Value = 5
$ws.Range("E649").Value = 2
$ws.Range("F649").Value = 1
$ws.Range("G649").Value = "Ischio"
$ws.Range("H649").Value = 0
$ws.Range("A650").Value = 45994
$ws.Range("B650").Value = "Hedi Nasri"
$ws.Range("C650").Value = 70
$ws.Range("D650").Value = 5
$ws.Range("E650").Value = 4
$ws.Range("F650").Value = 3
$ws.Range("G650").Value = "Ischio"
$ws.Range("H650").Value = 3
$ws.Range("A651").Value = 45994
$ws.Range("B651").Value = "Naim Dhib"
$ws.Range("C651").Value = 70
$ws.Range("D651").Value = 6
$ws.Range("E651").Value = 7
$ws.Range("F651").Value = 0
$ws.Range("H651").Value = 5
$ws.Range("A652").Value = 45994
$ws.Range("B652").Value = "Sofiane Belle"
$ws.Range("C652").Value = 70
$ws.Range("D652").Value = 6
$ws.Range("E652").Value = 4
$ws.Range("F652").Value = 2
$ws.Range("G652").Value = "Grnou"
$ws.Range("H652").Value = 6
$ws.Range("A653").Value = 45994
$ws.Range("B653").Value = "Maé Clavel"
$ws.Range("C653").Value = 70
$ws.Range("D653").Value = 6
$ws.Range("E653").Value = 6
$ws.Range("F653").Value = 6
$ws.Range("G653").Value = "Ischio"
$ws.Range("H653").Value = 6
$ws.Range("A654").Value = 45994
$ws.Range("B654").Value = "Kamal Bafounta"
$ws.Range("C654").Value = 70
$ws.Range("D654").Value = 7
$ws.Range("E654").Value = 7
$ws.Range("F654").Value = 3
$ws.Range("G654").Value = "Genou"
$ws.Range("H654").Value = 4


# --- Step 3: Formulas for column I (Charge = Volume * Intensite) ---
# Rows 639:643 are entered as one fill, rows 644:654 as another - matching the
# two shared-formula blocks visible in the source edit.
$ws.Range("I639:I643").Formula = "=C639*D639"
$ws.Range("I644:I654").Formula = "=C644*D644"

# --- Step 4: Update sheet view state to reflect where the user ended up editing ---
$win = $excel.ActiveWindow
$win.TopLeftCell = $ws.Range("A625")
$ws.Range("L642").Select()

